$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "neil"
$ws.Range("B3").Value = "farmer"
$ws.Range("C3").Value = "'1234"
$ws.Range("C3").Style = "Normal"

$ws.Range("A4").Value = "admin"
$ws.Range("B4").Value = "admin"
$ws.Range("C4").Value = "admin"
